# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (H,I,J,K,L,M,N) on several leve rows
# across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets with newly fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2171.2
$ws.Range("J40").Value = 2324.5
$ws.Range("L40").Value = 2324.5
$ws.Range("N40").Value = -2674.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 84995
$ws.Range("J87").Value = 84995
$ws.Range("L87").Value = 84995
$ws.Range("N87").Value = -87491

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 84995
$ws.Range("J90").Value = 84995
$ws.Range("L90").Value = 254985
$ws.Range("N90").Value = -267465

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 4104.7
$ws.Range("I135").Value = 4104.7
$ws.Range("K135").Value = 36942.3
$ws.Range("M135").Value = -34407.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2976.75
$ws.Range("I141").Value = 3139.5
$ws.Range("K141").Value = 9418.5
$ws.Range("M141").Value = -4238.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7037.9443
$ws.Range("I32").Value = 6193.273
$ws.Range("K32").Value = 6193.273
$ws.Range("M32").Value = -5906.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 31124.285
$ws.Range("J37").Value = 33328.332
$ws.Range("L37").Value = 33328.332
$ws.Range("N37").Value = -33874.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 64995
$ws.Range("J44").Value = 64995
$ws.Range("L44").Value = 64995
$ws.Range("N44").Value = -65971

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 75707.57000000001
$ws.Range("J55").Value = 76661.664
$ws.Range("L55").Value = 76661.664
$ws.Range("N55").Value = -77291.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10054.158
$ws.Range("I61").Value = 8162.4243
$ws.Range("K61").Value = 8162.4243
$ws.Range("M61").Value = -7950.4243

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 70369.25
$ws.Range("J80").Value = 72850.57000000001
$ws.Range("L80").Value = 72850.57000000001
$ws.Range("N80").Value = -74846.57000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 70369.25
$ws.Range("J83").Value = 72850.57000000001
$ws.Range("L83").Value = 218551.71
$ws.Range("N83").Value = -228535.71

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3667.9092
$ws.Range("J132").Value = 10432.833
$ws.Range("L132").Value = 31298.499
$ws.Range("N132").Value = -36358.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 10054.158
$ws.Range("I136").Value = 8162.4243
$ws.Range("K136").Value = 24487.2729
$ws.Range("M136").Value = -21937.2729

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 96661.664
$ws.Range("J35").Value = 96661.664
$ws.Range("L35").Value = 96661.664
$ws.Range("N35").Value = -97281.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 61338.668
$ws.Range("I82").Value = 13222
$ws.Range("J82").Value = 95707.71000000001
$ws.Range("K82").Value = 13222
$ws.Range("L82").Value = 95707.71000000001
$ws.Range("M82").Value = -12839
$ws.Range("N82").Value = -96473.71000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 61338.668
$ws.Range("I85").Value = 13222
$ws.Range("J85").Value = 95707.71000000001
$ws.Range("K85").Value = 13222
$ws.Range("L85").Value = 95707.71000000001
$ws.Range("M85").Value = -11896
$ws.Range("N85").Value = -98359.71000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 43328.332
$ws.Range("J41").Value = 43328.332
$ws.Range("L41").Value = 43328.332
$ws.Range("N41").Value = -44184.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 47995
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 45995.715
$ws.Range("I51").Value = 42000
$ws.Range("J51").Value = 46661.668
$ws.Range("K51").Value = 42000
$ws.Range("L51").Value = 46661.668
$ws.Range("M51").Value = -41264
$ws.Range("N51").Value = -48133.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5583.522
$ws.Range("I58").Value = 3221.4
$ws.Range("K58").Value = 3221.4
$ws.Range("M58").Value = -3018.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 45995.715
$ws.Range("I61").Value = 42000
$ws.Range("J61").Value = 46661.668
$ws.Range("K61").Value = 42000
$ws.Range("L61").Value = 46661.668
$ws.Range("M61").Value = -41652
$ws.Range("N61").Value = -47357.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 53686.1
$ws.Range("J74").Value = 53686.1
$ws.Range("L74").Value = 53686.1
$ws.Range("N74").Value = -55434.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 53686.1
$ws.Range("J77").Value = 53686.1
$ws.Range("L77").Value = 161058.3
$ws.Range("N77").Value = -169794.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 25621.363
$ws.Range("I132").Value = 14481.914
$ws.Range("J132").Value = 45115.4
$ws.Range("K132").Value = 43445.742
$ws.Range("L132").Value = 135346.2
$ws.Range("M132").Value = -40915.742
$ws.Range("N132").Value = -140406.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5486.2856
$ws.Range("I134").Value = 4197.346
$ws.Range("K134").Value = 12592.038
$ws.Range("M134").Value = -10057.038

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5583.522
$ws.Range("I136").Value = 3221.4
$ws.Range("K136").Value = 9664.200000000001
$ws.Range("M136").Value = -7114.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 35599.4
$ws.Range("J33").Value = 36999.75
$ws.Range("L33").Value = 36999.75
$ws.Range("N33").Value = -37503.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5751.5
$ws.Range("I126").Value = 5572.7144
$ws.Range("J126").Value = 6001.8
$ws.Range("K126").Value = 16718.1432
$ws.Range("L126").Value = 18005.4
$ws.Range("M126").Value = -14248.1432
$ws.Range("N126").Value = -22945.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 17599.158
$ws.Range("I132").Value = 9915.959999999999
$ws.Range("J132").Value = 32374.54
$ws.Range("K132").Value = 29747.88
$ws.Range("L132").Value = 97123.62
$ws.Range("M132").Value = -27217.88
$ws.Range("N132").Value = -102183.62

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 74996.28999999999
$ws.Range("J135").Value = 74996.28999999999
$ws.Range("L135").Value = 74996.28999999999
$ws.Range("N135").Value = -85136.28999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 31999.5
$ws.Range("J98").Value = 31999.5
$ws.Range("L98").Value = 31999.5
$ws.Range("N98").Value = -37989.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2589.4075
$ws.Range("I122").Value = 1953.421
$ws.Range("K122").Value = 5860.263
$ws.Range("M122").Value = -3410.263

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 391613.53
$ws.Range("I122").Value = 534134.3
$ws.Range("K122").Value = 1602402.9
$ws.Range("M122").Value = -1599952.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7410350.5
$ws.Range("I136").Value = 10528477
$ws.Range("K136").Value = 31585431
$ws.Range("M136").Value = -31582881
